$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy formatting (bold font, border, centered alignment) from H1 header cell
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows 2-11 for columns I (I0) and J (IF)
$data = @(
    @(6, 8),
    @(3, 4),
    @(11, 12),
    @(4, 8),
    @(1, 4),
    @(5, 8),
    @(4, 9),
    @(4, 6),
    @(5, 6),
    @(1, 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
